# Weekly update: insert a new observation as row 3 (pushing the existing
# rows 3-29 down to 4-30) for "Vega Monumental Concepción - Poroto granado".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 3..29 down to 4..30, leaving a blank row 3 for the new record.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the latest weekly observation.
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44602
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112030
$ws.Range("G3").Value = "Poroto granado"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 130
$ws.Range("K3").Value = 20000
$ws.Range("L3").Value = 21000
$ws.Range("M3").Value = 20385
$ws.Range("N3").Value = "`$/saco 25 kilos"
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 815
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
